$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Category" column (F) has blank cells for several rows (the source
# data had no category for those channels). Fill every blank cell in
# column F, within the used data range, with "Mixed" so the category
# filter/slicer works without blanks.

$lastRow = $ws.Cells.SpecialCells(11).Row   # xlCellTypeLastCell = 11

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)  # column F
    $val = $cell.Value()
    if ($null -eq $val -or $val -eq "") {
        $cell.Value = "Mixed"
    }
}

# Reflect the author's final UI selection: column F selected (whole column).
$ws.Range("F1:F1048576").Select()
